$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cage-seq")

# Insert three new columns before column B ("Assay Type") to make room for
# "Experiment Alias", "Project" and "Secondary Project" header cells.
# Column B currently holds "Assay Type"; after insertion we need:
#   B: Experiment Alias, C: Project, D: Secondary Project, E: Assay Type (shifted)
$ws.Range("B1:D1").EntireColumn.Insert()

$ws.Range("B1").Value = "Experiment Alias"
$ws.Range("C1").Value = "Project"
$ws.Range("D1").Value = "Secondary Project"

# Re-apply best-fit-style column widths across the whole used range (A:AI),
# matching the widths Excel computed when it auto-fit the header row after
# the new columns were added.
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.330729166666666
$ws.Columns.Item(3).ColumnWidth = 5.998697916666667
$ws.Columns.Item(4).ColumnWidth = 14.998697916666666
$ws.Columns.Item(5).ColumnWidth = 9.666666666666666
$ws.Columns.Item(6).ColumnWidth = 13.330729166666666
$ws.Columns.Item(7).ColumnWidth = 22.830729166666668
$ws.Columns.Item(8).ColumnWidth = 27.498697916666668
$ws.Columns.Item(9).ColumnWidth = 3.8307291666666665
$ws.Columns.Item(10).ColumnWidth = 18.666666666666668
$ws.Columns.Item(11).ColumnWidth = 15.998697916666666
$ws.Columns.Item(12).ColumnWidth = 23.666666666666668
$ws.Columns.Item(13).ColumnWidth = 32.498697916666664
$ws.Columns.Item(14).ColumnWidth = 3.8307291666666665
$ws.Columns.Item(15).ColumnWidth = 31.166666666666668
$ws.Columns.Item(16).ColumnWidth = 3.8307291666666665
$ws.Columns.Item(17).ColumnWidth = 20.830729166666668
$ws.Columns.Item(18).ColumnWidth = 3.8307291666666665
$ws.Columns.Item(19).ColumnWidth = 17.166666666666668
$ws.Columns.Item(20).ColumnWidth = 25.998697916666668
$ws.Columns.Item(21).ColumnWidth = 3.8307291666666665
$ws.Columns.Item(22).ColumnWidth = 24.498697916666668
$ws.Columns.Item(23).ColumnWidth = 3.8307291666666665
$ws.Columns.Item(24).ColumnWidth = 14.166666666666666
$ws.Columns.Item(25).ColumnWidth = 3.8307291666666665
$ws.Columns.Item(26).ColumnWidth = 15.830729166666666
$ws.Columns.Item(27).ColumnWidth = 12.998697916666666
$ws.Columns.Item(28).ColumnWidth = 11.998697916666666
$ws.Columns.Item(29).ColumnWidth = 23.330729166666668
$ws.Columns.Item(30).ColumnWidth = 22.666666666666668
$ws.Columns.Item(31).ColumnWidth = 18.830729166666668
$ws.Columns.Item(32).ColumnWidth = 30.998697916666668
$ws.Columns.Item(33).ColumnWidth = 20.998697916666668
$ws.Columns.Item(34).ColumnWidth = 20.998697916666668
$ws.Columns.Item(35).ColumnWidth = 18.830729166666668

$wb.Save()
